$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

# Update header row values (row 1)
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Update row 2 values; C2 and E2 get cleared (removed) entirely
$ws.Range("B2").Value = 21.45912128837421
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 30.2044515876383
$ws.Range("E2").ClearContents()

# Update row 3 values
$ws.Range("B3").Value = 18.257473325937074
$ws.Range("C3").Value = -7.7900079309787529
$ws.Range("D3").Value = 25.278493384463228
$ws.Range("E3").Value = -6.4305823250474115

# Update selection to match target sqref B1:E3
$ws.Range("B1:E3").Select()
